$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2409.2727
$ws.Range("I6").Value = 2625
$ws.Range("J6").Value = 2286
$ws.Range("K6").Value = 7875
$ws.Range("L6").Value = 6858
$ws.Range("M6").Value = -7763
$ws.Range("N6").Value = -7082

$ws.Range("H17").Value = 585212.94
$ws.Range("J17").Value = 585212.94
$ws.Range("L17").Value = 1755638.82
$ws.Range("N17").Value = -1755974.82

$ws.Range("H86").Value = 9093772
$ws.Range("I86").Value = 16667883
$ws.Range("J86").Value = 4838.8
$ws.Range("K86").Value = 16667883
$ws.Range("L86").Value = 4838.8
$ws.Range("M86").Value = -16666760
$ws.Range("N86").Value = -7084.8

$ws.Range("H89").Value = 9093772
$ws.Range("I89").Value = 16667883
$ws.Range("J89").Value = 4838.8
$ws.Range("K89").Value = 83339415
$ws.Range("L89").Value = 24194
$ws.Range("M89").Value = -83333799
$ws.Range("N89").Value = -35426

$ws.Range("H127").Value = 39339.46
$ws.Range("I127").Value = 167057.67
$ws.Range("J127").Value = 1024
$ws.Range("K127").Value = 501173.01
$ws.Range("L127").Value = 3072
$ws.Range("M127").Value = -496213.01
$ws.Range("N127").Value = -12992

$ws.Range("H129").Value = 530.2105
$ws.Range("I129").Value = 406
$ws.Range("J129").Value = 878
$ws.Range("K129").Value = 1218
$ws.Range("L129").Value = 2634
$ws.Range("M129").Value = 3782
$ws.Range("N129").Value = -12634

$ws.Range("H137").Value = 36808.965
$ws.Range("I137").Value = 67609
$ws.Range("J137").Value = 1270.4615
$ws.Range("K137").Value = 202827
$ws.Range("L137").Value = 3811.3845
$ws.Range("M137").Value = -200277
$ws.Range("N137").Value = -8911.3845

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1029.1428
$ws.Range("I61").Value = 882.8
$ws.Range("J61").Value = 1395
$ws.Range("K61").Value = 882.8
$ws.Range("L61").Value = 1395
$ws.Range("M61").Value = -670.8
$ws.Range("N61").Value = -1819

$ws.Range("H136").Value = 1029.1428
$ws.Range("I136").Value = 882.8
$ws.Range("J136").Value = 1395
$ws.Range("K136").Value = 2648.4
$ws.Range("L136").Value = 4185
$ws.Range("M136").Value = -98.39999999999964
$ws.Range("N136").Value = -9285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35919.41
$ws.Range("I31").Value = 38542.406
$ws.Range("J31").Value = 23928.572
$ws.Range("K31").Value = 38542.406
$ws.Range("L31").Value = 23928.572
$ws.Range("M31").Value = -38247.406
$ws.Range("N31").Value = -24518.572

$ws.Range("H34").Value = 35919.41
$ws.Range("I34").Value = 38542.406
$ws.Range("J34").Value = 23928.572
$ws.Range("K34").Value = 38542.406
$ws.Range("L34").Value = 23928.572
$ws.Range("M34").Value = -38340.406
$ws.Range("N34").Value = -24332.572

$ws.Range("H62").Value = 3054.2
$ws.Range("I62").Value = 3004.6667
$ws.Range("K62").Value = 3004.6667
$ws.Range("M62").Value = -2380.6667

$ws.Range("H65").Value = 3054.2
$ws.Range("I65").Value = 3004.6667
$ws.Range("K65").Value = 15023.3335
$ws.Range("M65").Value = -11903.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 2861.5386
$ws.Range("J33").Value = 7732.75
$ws.Range("L33").Value = 46396.5
$ws.Range("N33").Value = -46962.5

$ws.Range("H44").Value = 187501680
$ws.Range("I44").Value = 500000260
$ws.Range("J44").Value = 83335500
$ws.Range("K44").Value = 1500000780
$ws.Range("L44").Value = 250006500
$ws.Range("M44").Value = -1500000382
$ws.Range("N44").Value = -250007296

$ws.Range("H46").Value = 2232.3845
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 2335.0833
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 7005.249899999999
$ws.Range("M46").Value = -2909
$ws.Range("N46").Value = -7187.249899999999

$ws.Range("H64").Value = 2168394.8
$ws.Range("I64").Value = 1470.6666
$ws.Range("J64").Value = 2759374
$ws.Range("K64").Value = 4411.9998
$ws.Range("L64").Value = 8278122
$ws.Range("M64").Value = -4141.9998
$ws.Range("N64").Value = -8278662

$ws.Range("H67").Value = 2168394.8
$ws.Range("I67").Value = 1470.6666
$ws.Range("J67").Value = 2759374
$ws.Range("K67").Value = 4411.9998
$ws.Range("L67").Value = 8278122
$ws.Range("M67").Value = -3475.9998
$ws.Range("N67").Value = -8279994

$ws.Range("H92").Value = 1333.3334
$ws.Range("J92").Value = 2000
$ws.Range("L92").Value = 6000
$ws.Range("N92").Value = -8496

$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("M98").ClearContents()

$ws.Range("H127").Value = 2137500
$ws.Range("J127").Value = 2137500
$ws.Range("L127").Value = 6412500
$ws.Range("N127").Value = -6422420

$ws.Range("H130").Value = 2833.3333
$ws.Range("I130").Value = 1966.6666
$ws.Range("J130").Value = 3700
$ws.Range("K130").Value = 5899.9998
$ws.Range("L130").Value = 11100
$ws.Range("M130").Value = -879.9997999999996
$ws.Range("N130").Value = -21140

$ws.Range("H131").Value = 22728134
$ws.Range("I131").Value = 676.6667
$ws.Range("J131").Value = 24039334
$ws.Range("K131").Value = 2030.0001
$ws.Range("L131").Value = 72118002
$ws.Range("M131").Value = 3009.9999
$ws.Range("N131").Value = -72128082

$ws.Range("H134").Value = 4120.524
$ws.Range("I134").Value = 1155.3636
$ws.Range("J134").Value = 7382.2
$ws.Range("K134").Value = 3466.0908
$ws.Range("L134").Value = 22146.6
$ws.Range("M134").Value = 1603.9092
$ws.Range("N134").Value = -32286.6

$ws.Range("H138").Value = 7410089
$ws.Range("I138").Value = 1256
$ws.Range("J138").Value = 11114505
$ws.Range("K138").Value = 3768
$ws.Range("L138").Value = 33343515
$ws.Range("M138").Value = 1372
$ws.Range("N138").Value = -33353795

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 16944.5
$ws.Range("J45").Value = 16944.5
$ws.Range("L45").Value = 16944.5
$ws.Range("N45").Value = -18062.5

$ws.Range("H70").Value = 2568633
$ws.Range("I70").Value = 3452565.2
$ws.Range("J70").Value = 5229.6
$ws.Range("K70").Value = 3452565.2
$ws.Range("L70").Value = 5229.6
$ws.Range("M70").Value = -3452295.2
$ws.Range("N70").Value = -5769.6

$ws.Range("H73").Value = 2568633
$ws.Range("I73").Value = 3452565.2
$ws.Range("J73").Value = 5229.6
$ws.Range("K73").Value = 3452565.2
$ws.Range("L73").Value = 5229.6
$ws.Range("M73").Value = -3451629.2
$ws.Range("N73").Value = -7101.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 46895
$ws.Range("J123").Value = 46895
$ws.Range("L123").Value = 46895
$ws.Range("N123").Value = -56695

$ws.Range("H136").Value = 372537.84
$ws.Range("I136").Value = 589869
$ws.Range("J136").Value = 3074.9
$ws.Range("K136").Value = 1769607
$ws.Range("L136").Value = 9224.7
$ws.Range("M136").Value = -1767057
$ws.Range("N136").Value = -14324.7
